$d = $word.ActiveDocument

# Locate the sentence that ends the list item; the two trailing
# "empty" runs (each holding only a manual line break <w:br/>) that
# follow it are what we need to remove.
$findRange = $d.Content
$found = $findRange.Find.Execute("This pair of numbers is the answer.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $afterText = $findRange.End
    $para = $findRange.Paragraphs(1)
    $paraEnd = $para.Range.End

    # $paraEnd - 1 stops just before the paragraph mark, so this range
    # covers only the trailing line breaks that belong to the two
    # break-only runs, leaving the paragraph mark itself untouched.
    if ($paraEnd - 1 -gt $afterText) {
        $trailing = $d.Range($afterText, $paraEnd - 1)
        $trailing.Delete()
    }
}
